$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (shared string index 224, cell A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 15:16"

# Swap Croacia / Luxemburgo entries: row 101 (was Luxemburgo) becomes Croacia,
# row 102 (was Croacia) becomes Luxemburgo, matching the shared-string reorder in the diff.
$ws.Range("A101").Value = "Croacia"
$ws.Range("A102").Value = "Luxemburgo"

# Update per-country statistic values
$ws.Range("B4").Value = 5800392
$ws.Range("C4").Value = 3665
$ws.Range("D4").Value = 3127438
$ws.Range("E4").Value = 2493706
$ws.Range("G4").Value = 48
$ws.Range("H4").Value = 179248

$ws.Range("B6").Value = 2985367
$ws.Range("C6").Value = 11999
$ws.Range("D6").Value = 2228248
$ws.Range("E6").Value = 701089
$ws.Range("G6").Value = 102
$ws.Range("H6").Value = 56030

$ws.Range("B17").Value = 306370
$ws.Range("C17").Value = 1184
$ws.Range("D17").Value = 278441
$ws.Range("E17").Value = 24310
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = 3619

$ws.Range("B24").Value = 201050
$ws.Range("C24").Value = 3965
$ws.Range("D24").Value = 143393
$ws.Range("E24").Value = 51304
$ws.Range("G24").Value = 70
$ws.Range("H24").Value = 6353

$ws.Range("B43").Value = 70285
$ws.Range("C43").Value = 174
$ws.Range("D43").Value = 68577
$ws.Range("E43").Value = 1071
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 637

$ws.Range("B46").Value = 66097
$ws.Range("C46").Value = 508
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 6200

$ws.Range("B60").Value = 38412
$ws.Range("C60").Value = 338
$ws.Range("E60").Value = 4158
$ws.Range("G60").Value = 5
$ws.Range("H60").Value = 265

$ws.Range("B68").Value = 31117
$ws.Range("C68").Value = 634
$ws.Range("D68").Value = 18350
$ws.Range("E68").Value = 12621
$ws.Range("G68").Value = 9
$ws.Range("H68").Value = 146

$ws.Range("B80").Value = 16239
$ws.Range("C80").Value = 112
$ws.Range("E80").Value = 1673
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 622

$ws.Range("B84").Value = 12850
$ws.Range("C84").Value = 161
$ws.Range("E84").Value = 4419
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 266

$ws.Range("B92").Value = 9257
$ws.Range("C92").Value = 8
$ws.Range("D92").Value = 8949
$ws.Range("E92").Value = 183

$ws.Range("B101").Value = 7900
$ws.Range("C101").Value = 306
$ws.Range("D101").Value = 5678
$ws.Range("E101").Value = 2052
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 170

$ws.Range("D102").Value = 6969
$ws.Range("E102").Value = 611
$ws.Range("H102").Value = 124

$ws.Range("B138").Value = 2058
$ws.Range("C138").Value = 8
$ws.Range("D138").Value = 1936
$ws.Range("E138").Value = 112

$ws.Range("B144").Value = 1577
$ws.Range("C144").Value = 31
$ws.Range("D144").Value = 901
$ws.Range("E144").Value = 666
